$wb = $excel.ActiveWorkbook

# --- Summary sheet: updated aggregate stats ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = 0.48
$summary.Range("B6").Value = 5
$summary.Range("B9").Value = 40

# --- All Trades sheet ---
$allTrades = $wb.Worksheets.Item("All Trades")

# Header row: "Exit Reason" and "Duration (min)" moved earlier, right after "Capital After"
$allTrades.Range("L1").Value = "Exit Reason"
$allTrades.Range("M1").Value = "Duration (min)"
$allTrades.Range("N1").Value = "Entry Slippage (bps)"
$allTrades.Range("O1").Value = "Exit Slippage (bps)"
$allTrades.Range("P1").Value = "Confidence"
$allTrades.Range("Q1").Value = "Entry Reason"

# Trade #5 (row 6) closes out
$allTrades.Range("G6").Value = 0.01
$allTrades.Range("H6").Value = "CLOSED"
$allTrades.Range("K6").Value = 100
$allTrades.Range("L6").Value = "early_exit"
$allTrades.Range("M6").Value = 0.1

# Existing rows 24-27: re-seat the columns that moved (Entry/Exit Slippage, Confidence, Entry Reason)
# column L held "Entry Slippage" before and is now blank ("Exit Reason", not populated on these closed/open legacy rows)
$allTrades.Range("L24").Value = ""
$allTrades.Range("M24").Value = 0
$allTrades.Range("N24").Value = 0
$allTrades.Range("O24").Value = 0
$allTrades.Range("P24").Value = 0.6
$allTrades.Range("Q24").Value = "Normal spread capture: 408 bps"

$allTrades.Range("L25").Value = ""
$allTrades.Range("M25").Value = 0
$allTrades.Range("N25").Value = 0
$allTrades.Range("O25").Value = 0
$allTrades.Range("P25").Value = 0.6
$allTrades.Range("Q25").Value = "Normal spread capture: 202 bps"

$allTrades.Range("L26").Value = ""
$allTrades.Range("M26").Value = 0
$allTrades.Range("N26").Value = 0
$allTrades.Range("O26").Value = 0
$allTrades.Range("P26").Value = 0.6
$allTrades.Range("Q26").Value = "Normal spread capture: 202 bps"

$allTrades.Range("L27").Value = ""
$allTrades.Range("M27").Value = 0
$allTrades.Range("N27").Value = 0
$allTrades.Range("O27").Value = 0
$allTrades.Range("P27").Value = 0.9
$allTrades.Range("Q27").Value = "Upward momentum: 1.020% over 5 samples"

# New trade rows logged to the master "All Trades" sheet (new column order)
$allTrades.Range("A28").Value = 27
$allTrades.Range("B28").Value = "'2026-02-18"
$allTrades.Range("C28").Value = "10:30:00"
$allTrades.Range("D28").Value = "momentum"
$allTrades.Range("E28").Value = "UP"
$allTrades.Range("F28").Value = 0.01
$allTrades.Range("H28").Value = "OPEN"
$allTrades.Range("I28").Value = 0
$allTrades.Range("J28").Value = 0
$allTrades.Range("K28").Value = 100
$allTrades.Range("M28").Value = 0
$allTrades.Range("N28").Value = 0
$allTrades.Range("O28").Value = 0
$allTrades.Range("P28").Value = 0.9
$allTrades.Range("Q28").Value = "Upward momentum: 1.020% over 9 samples"

$allTrades.Range("A29").Value = 28
$allTrades.Range("B29").Value = "'2026-02-18"
$allTrades.Range("C29").Value = "10:30:01"
$allTrades.Range("D29").Value = "MarketMaking"
$allTrades.Range("E29").Value = "UP"
$allTrades.Range("F29").Value = 0.01
$allTrades.Range("H29").Value = "OPEN"
$allTrades.Range("I29").Value = 0
$allTrades.Range("J29").Value = 0
$allTrades.Range("K29").Value = 100
$allTrades.Range("M29").Value = 0
$allTrades.Range("N29").Value = 0
$allTrades.Range("O29").Value = 0
$allTrades.Range("P29").Value = 0.6
$allTrades.Range("Q29").Value = "Normal spread capture: 202 bps"

# --- momentum strategy sheet: new open trade logged (original column order) ---
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A3").Value = 27
$momentum.Range("B3").Value = "'2026-02-18"
$momentum.Range("C3").Value = "10:30:00"
$momentum.Range("D3").Value = "momentum"
$momentum.Range("E3").Value = "UP"
$momentum.Range("F3").Value = 0.01
$momentum.Range("H3").Value = "OPEN"
$momentum.Range("I3").Value = 0
$momentum.Range("J3").Value = 0
$momentum.Range("K3").Value = 100
$momentum.Range("L3").Value = 0
$momentum.Range("M3").Value = 0
$momentum.Range("N3").Value = 0.9
$momentum.Range("O3").Value = "Upward momentum: 1.020% over 9 samples"
$momentum.Range("Q3").Value = 0

# --- MarketMaking strategy sheet: new open trade logged (original column order) ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("A5").Value = 28
$marketMaking.Range("B5").Value = "'2026-02-18"
$marketMaking.Range("C5").Value = "10:30:01"
$marketMaking.Range("D5").Value = "MarketMaking"
$marketMaking.Range("E5").Value = "UP"
$marketMaking.Range("F5").Value = 0.01
$marketMaking.Range("H5").Value = "OPEN"
$marketMaking.Range("I5").Value = 0
$marketMaking.Range("J5").Value = 0
$marketMaking.Range("K5").Value = 100
$marketMaking.Range("L5").Value = 0
$marketMaking.Range("M5").Value = 0
$marketMaking.Range("N5").Value = 0.6
$marketMaking.Range("O5").Value = "Normal spread capture: 202 bps"
$marketMaking.Range("Q5").Value = 0

